# Update the cryptos list (Price and Volume(1h) columns) with freshly
# scraped values, as produced by the GitHub Actions scheduled job.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: spreadsheet row number, new Price (column D, $null = unchanged),
# new Volume(1h) text (column E).
$updates = @(
    @{ Row = 2;  D = "64.057.58";  E = "  -0.35%  " },
    @{ Row = 3;  D = "3.471.54";   E = "  -0.88%  " },
    @{ Row = 4;  D = "1.00";       E = "  +0.04%  " },
    @{ Row = 5;  D = "583.69";     E = "  -0.49%  " },
    @{ Row = 6;  D = "130.85";     E = "  -2.58%  " },
    @{ Row = 7;  D = $null;        E = "  +0.04%  " },
    @{ Row = 8;  D = $null;        E = "  -1.08%  " },
    @{ Row = 9;  D = "7.59";       E = "  +4.61%  " },
    @{ Row = 10; D = $null;        E = "  -1.61%  " },
    @{ Row = 11; D = "0.387";      E = "  -0.05%  " },
    @{ Row = 12; D = "4.071.05";   E = "  -0.67%  " },
    @{ Row = 13; D = $null;        E = "  -0.11%  " },
    @{ Row = 14; D = $null;        E = "  -3.20%  " },
    @{ Row = 15; D = "3.471.94";   E = "  -0.87%  " },
    @{ Row = 16; D = "64.051.62";  E = "  -0.39%  " },
    @{ Row = 17; D = "24.24";      E = "  -6.66%  " },
    @{ Row = 18; D = "9.94";       E = "  +0.11%  " },
    @{ Row = 19; D = "5.67";       E = "  -1.43%  " },
    @{ Row = 20; D = "13.43";      E = "  -1.84%  " },
    @{ Row = 21; D = "383.84";     E = "  -2.40%  " },
    @{ Row = 22; D = "0.568";      E = "  -0.81%  " },
    @{ Row = 23; D = "3.613.87";   E = "  -0.80%  " },
    @{ Row = 24; D = "74.89";      E = "  +0.82%  " },
    @{ Row = 25; D = $null;        E = "  +0.07%  " },
    @{ Row = 26; D = "5.61";       E = "  +0.45%  " },
    @{ Row = 27; D = "0.0000111";  E = "  -3.58%  " },
    @{ Row = 28; D = $null;        E = "  +0.02%  " },
    @{ Row = 29; D = "2.22";       E = "  -0.65%  " },
    @{ Row = 30; D = "7.04";       E = "  -4.92%  " },
    @{ Row = 31; D = $null;        E = "  -4.31%  " },
    @{ Row = 32; D = "7.89";       E = "  -4.76%  " },
    @{ Row = 33; D = "3.502.63";   E = "  -0.58%  " },
    @{ Row = 34; D = $null;        E = "  +0.94%  " },
    @{ Row = 36; D = "22.84";      E = "  -2.69%  " },
    @{ Row = 37; D = "5.19";       E = "  +0.15%  " },
    @{ Row = 38; D = "6.74";       E = "  -2.70%  " },
    @{ Row = 39; D = "1.49";       E = "  -4.57%  " },
    @{ Row = 40; D = "161.98";     E = "  -1.20%  " },
    @{ Row = 41; D = "0.0775";     E = "  -1.11%  " },
    @{ Row = 42; D = "0.796";      E = "  -1.35%  " },
    @{ Row = 43; D = "1.00";       E = "  +0.15%  " },
    @{ Row = 44; D = "41.33";      E = "  -1.08%  " },
    @{ Row = 45; D = "4.28";       E = "  -3.39%  " },
    @{ Row = 46; D = "1.61";       E = "  -2.58%  " },
    @{ Row = 47; D = $null;        E = "  -3.43%  " },
    @{ Row = 48; D = "23.32";      E = "  -7.42%  " },
    @{ Row = 49; D = "6.68";       E = "  -1.54%  " },
    @{ Row = 50; D = "0.900";      E = "  +0.07%  " },
    @{ Row = 51; D = "2.328.82";   E = "  -5.37%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $ws.Range("D$($u.Row)").Value = $u.D
    }
    $ws.Range("E$($u.Row)").Value = $u.E
}
